$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 177, pushing existing rows 177-215 down to 178-216.
$ws.Rows.Item(177).Insert()

# Populate the newly inserted row 177 with the new weekly record.
$ws.Cells.Item(177, 1).Value = 3
$ws.Cells.Item(177, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(177, 3).Value = "Coquimbo"
$ws.Cells.Item(177, 4).Value = 44943
$ws.Cells.Item(177, 5).Value = 5
$ws.Cells.Item(177, 6).Value = 100112052
$ws.Cells.Item(177, 7).Value = "Albahaca"
$ws.Cells.Item(177, 8).Value = "Sin especificar"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 80
$ws.Cells.Item(177, 11).Value = 4000
$ws.Cells.Item(177, 12).Value = 4000
$ws.Cells.Item(177, 13).Value = 4000
$ws.Cells.Item(177, 14).Value = "`$/docena de matas"
$ws.Cells.Item(177, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(177, 16).Value = 667
$ws.Cells.Item(177, 17).Value = 6
$ws.Cells.Item(177, 18).Value = "Hortaliza"
